# vessel_data.xlsx — "added new vessel_data and vessel_name"
#
# 1) Row 101 (TYANA) was originally exported with every numeric-looking field
#    stored as text (and missing fields stored as the literal text "None").
#    Fix the numeric columns to be real numbers, and blank out the "None"
#    placeholder cells.
# 2) Append a brand-new row 102 for vessel "MP ULTRAMAX 1", using the same
#    (text-based) export shape as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 101 — TYANA: retype numeric columns, blank the "None" placeholders
# ---------------------------------------------------------------------------
$ws.Cells.Item(101, 3).Value  = 9418468     # C101 IMO Number
$ws.Cells.Item(101, 4).Value  = ""          # D101 Launch Date ("None" -> blank)
$ws.Cells.Item(101, 5).Value  = 82158       # E101 Design Deadweight
$ws.Cells.Item(101, 6).Value  = 42965       # F101 Gross Tonnage
$ws.Cells.Item(101, 8).Value  = 94590       # H101 Displacement (tonnes)
$ws.Cells.Item(101, 9).Value  = ""          # I101 Design Speed Ahead ("None" -> blank)
$ws.Cells.Item(101, 10).Value = 228.99      # J101 Length Overall (LOA)
$ws.Cells.Item(101, 11).Value = 5649        # K101 Calculated Freeboard (mm)
$ws.Cells.Item(101, 12).Value = ""          # L101 Breadth Overall ("None" -> blank)
$ws.Cells.Item(101, 13).Value = ""          # M101 Depth Overall ("None" -> blank)
$ws.Cells.Item(101, 14).Value = 222.54      # N101 Length Between Perpendicular (LPP)
$ws.Cells.Item(101, 15).Value = 14.434      # O101 Design Draft
$ws.Cells.Item(101, 16).Value = ""          # P101 Draft Molded ("None" -> blank)
$ws.Cells.Item(101, 17).Value = ""          # Q101 Draft Scantling ("None" -> blank)
$ws.Cells.Item(101, 19).Value = ""          # S101 Auxiliary Engine Rated Power ("None" -> blank)
$ws.Cells.Item(101, 20).Value = ""          # T101 Auxiliary Engine Manufacturer ("None" -> blank)
$ws.Cells.Item(101, 24).Value = 9710        # X101 Main Engine Rated Power

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (even if it looks numeric) without leaving a
# permanent "Text" number-format on the cell — match the plain/unstyled
# look of the rest of the sheet.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Row 102 — new vessel: MP ULTRAMAX 1 (Republic of Singapore)
# ---------------------------------------------------------------------------
$ws.Cells.Item(102, 1).Value = "MP ULTRAMAX 1"          # A102 Vessel Name
$ws.Cells.Item(102, 2).Value = "Republic of Singapore"  # B102 Flag

Set-TextValue $ws.Cells.Item(102, 3)  "9703590"    # C102 IMO Number
$ws.Cells.Item(102, 4).Value = "None"               # D102 Launch Date
Set-TextValue $ws.Cells.Item(102, 5)  "63339.85"   # E102 Design Deadweight
Set-TextValue $ws.Cells.Item(102, 6)  "36286"      # F102 Gross Tonnage
$ws.Cells.Item(102, 7).Value = ""                   # G102 Block Coefficient Class (Cb)
Set-TextValue $ws.Cells.Item(102, 8)  "75196.6"    # H102 Displacement (tonnes)
Set-TextValue $ws.Cells.Item(102, 9)  "14.4"       # I102 Design Speed Ahead
Set-TextValue $ws.Cells.Item(102, 10) "199.9"      # J102 Length Overall (LOA)
Set-TextValue $ws.Cells.Item(102, 11) "5228"       # K102 Calculated Freeboard (mm)
$ws.Cells.Item(102, 12).Value = "None"              # L102 Breadth Overall
$ws.Cells.Item(102, 13).Value = "None"              # M102 Depth Overall
Set-TextValue $ws.Cells.Item(102, 14) "194.5"      # N102 Length Between Perpendicular (LPP)
Set-TextValue $ws.Cells.Item(102, 15) "11.3"       # O102 Design Draft
Set-TextValue $ws.Cells.Item(102, 16) "13.3"       # P102 Draft Molded
$ws.Cells.Item(102, 17).Value = "None"              # Q102 Draft Scantling
$ws.Cells.Item(102, 18).Value = 3                   # R102 Generator Number
Set-TextValue $ws.Cells.Item(102, 19) "720"        # S102 Auxiliary Engine Rated Power
$ws.Cells.Item(102, 20).Value = "None"              # T102 Auxiliary Engine Manufacturer
$ws.Cells.Item(102, 21).Value = ""                  # U102 Auxiliary Engine Model Number
$ws.Cells.Item(102, 22).Value = "None"              # V102 Main Engine Manufacturer
$ws.Cells.Item(102, 23).Value = ""                  # W102 Main Engine Model Number
$ws.Cells.Item(102, 24).Value = "None"              # X102 Main Engine Rated Power
